$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-10 from 2023-09-01 (45170)
# to 2023-09-05 (45174), reflecting an automatic data refresh.
$ws.Range("C2:C10").Value = 45174
